# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
$wb = $excel.ActiveWorkbook

# OFF sheet: row 2 (A2 = "H") — update Short Att, Short Comp, Deep Att, Deep Comp
$wsOFF = $wb.Worksheets.Item("OFF")
$wsOFF.Range("B2").Value = 418
$wsOFF.Range("C2").Value = 305
$wsOFF.Range("D2").Value = 105
$wsOFF.Range("E2").Value = 48

# DEF sheet: row 2 (A2 = "H") — update Short Att, Short Comp, Deep Att, Deep Comp, Short Int, Deep Int
$wsDEF = $wb.Worksheets.Item("DEF")
$wsDEF.Range("B2").Value = 629
$wsDEF.Range("C2").Value = 457
$wsDEF.Range("D2").Value = 119
$wsDEF.Range("E2").Value = 56
$wsDEF.Range("F2").Value = 5
$wsDEF.Range("G2").Value = 8
